$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet's display name (the <sheet name="..."> in workbook.xml).
$ws.Name = "GammaFiber2F"

# 2) Add new row 16, mirroring row 15's formatting (bold/bordered index cell in column A,
#    normal cells elsewhere), then fill in its values:
#      A16 = 14 (continues the running index)
#      B16 = "HexGrid-60degTilt5degRes" (same category repeated from row 15)
#      C16:M16 = 1
$ws.Range("A15:M15").Copy() | Out-Null
$ws.Range("A16:M16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1
